# Update cryptos list (prices / 1h volume %) as scraped by the GitHub Actions job.
# Cells in column D whose new value would otherwise be auto-parsed by Excel as a
# number (losing formatting such as trailing zeros, e.g. "12.00" -> 12) are
# forced to remain text via NumberFormat "@" before the value is written.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.157.20'
$ws.Range("E2").Value = '  -1.92%  '

$ws.Range("D3").Value = '2.589.59'
$ws.Range("E3").Value = '  -2.24%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.96'
$ws.Range("E5").Value = '  -3.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.40'
$ws.Range("E6").Value = '  -2.71%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.581'
$ws.Range("E8").Value = '  -1.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.108'
$ws.Range("E9").Value = '  -1.73%  '

$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.383'
$ws.Range("E11").Value = '  -2.14%  '

$ws.Range("E12").Value = '  -0.80%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.34'
$ws.Range("E13").Value = '  -2.10%  '

$ws.Range("D14").Value = '3.051.19'
$ws.Range("E14").Value = '  -2.53%  '

$ws.Range("D15").Value = '62.988.29'
$ws.Range("E15").Value = '  -1.95%  '

$ws.Range("E16").Value = '  +1.35%  '

$ws.Range("D17").Value = '2.609.58'
$ws.Range("E17").Value = '  -1.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.00'
$ws.Range("E18").Value = '  -1.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.62'
$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '342.77'
$ws.Range("E20").Value = '  -2.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.76'
$ws.Range("E21").Value = '  -2.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.22'
$ws.Range("E23").Value = '  -0.92%  '

$ws.Range("E24").Value = '  -2.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.10'
$ws.Range("E25").Value = '  -3.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.62'
$ws.Range("E26").Value = '  -5.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '558.45'
$ws.Range("E27").Value = '  +2.61%  '

$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.08'
$ws.Range("E28").Value = '  -2.65%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.162'
$ws.Range("E29").Value = '  -3.15%  '

$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.00'
$ws.Range("E31").Value = '  -3.60%  '

$ws.Range("D32").Value = '0.0₃0833'
$ws.Range("E32").Value = '  -3.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.73'
$ws.Range("E33").Value = '  -1.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.29'
$ws.Range("E34").Value = '  -1.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '166.02'
$ws.Range("E35").Value = '  -1.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.408'
$ws.Range("E36").Value = '  -0.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.22'
$ws.Range("E38").Value = '  -2.13%  '

$ws.Range("E39").Value = '  -5.52%  '

$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '165.76'
$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.94'
$ws.Range("E42").Value = '  +0.31%  '

$ws.Range("E43").Value = '  +3.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0577'
$ws.Range("E44").Value = '  -0.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.623'
$ws.Range("E45").Value = '  -1.22%  '

$ws.Range("E46").Value = '  +0.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0245'
$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0954'
$ws.Range("E48").Value = '  -1.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.86'
$ws.Range("E49").Value = '  -2.60%  '

$ws.Range("D50").Value = '0.0₆0223'
$ws.Range("E50").Value = '  +10.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.177'
$ws.Range("E51").Value = '  -5.80%  '
